$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2026 -Month 2 -Day 9 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($r = 2; $r -le 232; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
